$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New MFH building rows (7-10), matching the existing SFH row data layout.
$rows = @(
    @{ A=6; B="MFH"; C=2002; D=2008; F=526.54188049999993;  G=545.92182505920346;  H=110.37063172236979; I=195.04737837333329; J=1976692.2205915591; K=3.1881206593694991; L=4.09; M=19.985595685714241; N=9.9927978428571205;  O=9.9927978428571205;  P=21000 },
    @{ A=7; B="MFH"; C=2002; D=2008; F=526.54188049999993;  G=311.63072655812732;  H=73.936508821704834; I=192.79152078666669; J=1825798.5110075211; K=3.3021848795015778; L=4.09; M=18.174420757741562; N=9.087210378870779;   O=9.087210378870779;   P=21000 },
    @{ A=8; B="MFH"; C=2002; D=2008; F=390.84052000000003;  G=715.49033468132347;  H=77.448621722548552; I=177.82566969333331; J=3482782.0897382898; K=3.5556844994939119; L=4.09; M=12.228544698146701; N=6.1142723490733522;  O=6.1142723490733522;  P=21000 },
    @{ A=9; B="MFH"; C=2002; D=2008; F=390.84052000000003;  G=363.4430590073618;   H=50.280800635531001; I=177.96626804666661; J=3020247.2967746542; K=3.6954541563476719; L=4.09; M=11.250767088751211; N=5.6253835443756044;  O=5.6253835443756044;  P=21000 }
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $r++
}

# Fill E7:E10 with one relative formula so Excel records it as a shared
# formula group (matches the authored file's <f t="shared" .../> layout).
$ws.Range("E7:E10").Formula = "=F7/42.5"

# E6 was previously the tail of shared group E3:E6 (si=0). Since the new
# rows start a fresh shared group at E7, re-stamp E6 with its own explicit
# formula so it no longer references the old shared group.
$ws.Range("E6").Formula = "=F6/42.5"

$ws.Range("D17").Select()
